$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at L, shifting L..Q to M..R
$ws.Columns("L").Insert()
$ws.Columns("L").ColumnWidth = 12.285714285714286

# Row 1
$ws.Range("N1").Value = 'These Canadian titans of sketch comedy continue to influence funny people more than two decades after their CBC/CBS/HBO show went off the air.'

# Row 2
$ws.Range("N2").Value = 'The 81-year-old entertainer turns banter with the audience into a whole show.'

# Row 3
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 'New York'
$ws.Range("C3").Value = 'Dark Universe'
$ws.Range("D3").Value = 'American Museum of Natural History'
$ws.Range("E3").Value = 10024
$ws.Range("F3").Value = 'Central Park West at 79th Street, Manhattan, NY 10024'
$ws.Range("G3").Value = '212-769-5100'
$ws.Range("H3").Value = ' amnh.org'
$ws.Range("J3").Value = 20141102
$ws.Range("N3").Value = 'The American Museum of Natural History presents its newest space show.'
$ws.Range("O3").Value = 'The American Museum of Natural History presents its newest space show, Dark Universe, which celebrates the discoveries that have led us to a greater level of knowledge about our universe, its history and our planet''s place in it. Expect breathtaking renderings of cosmic phenomena and spectacular scenes that will make you feel like you''ve traveled into space. Astrophysicist Neil deGrasse Tyson narrates. '
$ws.Range("P3").Value = 'http://www.nycgo.com/images/460x285/DarkUniverse_AMNH_V1_460x285.jpg'

# Row 4
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 'New York'
$ws.Range("C4").Value = 'Mapping Brooklyn'
$ws.Range("D4").Value = 'BRIC House'
$ws.Range("E4").Value = 11217
$ws.Range("F4").Value = '647 Fulton St., Brooklyn, NY 11217'
$ws.Range("G4").Value = '718-683-5600'
$ws.Range("H4").Value = ' bricartsmedia.org'
$ws.Range("J4").Value = 20150226
$ws.Range("K4").Value = 20150926
$ws.Range("N4").Value = 'Mapping Brooklyn features both historic maps and contemporary works of art that make use of mapping and cartography. A'
$ws.Range("O4").Value = 'Mapping Brooklyn features both historic maps and contemporary works of art that make use of mapping and cartography. Among the highlights of the exhibition, which will be on view at BRIC (February 26-May 3) and the Brooklyn Historical Society (February 26-September 6), are a colorful pictorial road map to the 1939 New York World’s Fair, a commercial edition of a Red Scare-era map depicting enclaves of suspected radical activity and a detailed map of one of Brooklyn’s earliest botanic gardens. These are complemented by contemporary works by artists who researched BHS'' map collection for their art including pieces by Justin Blinder, Christine Gedeon, Nick Vaughan & Jake Margolin and Sarah Williams.'
$ws.Range("P4").Value = 'http://www.nycgo.com/images/460x285/MappingBrooklyn_V1_460x285.jpg'

# Row 5
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 'New York'
$ws.Range("C5").Value = 'Christopher Wool'
$ws.Range("D5").Value = 'Luhring Augustine'
$ws.Range("E5").Value = 10011
$ws.Range("F5").Value = '531 W 24th St, Manhattan, NY 10011'
$ws.Range("G5").Value = '212-206-9100'
$ws.Range("H5").Value = ' luhringaugustine.com'
$ws.Range("J5").Value = 20150502
$ws.Range("K5").Value = 20150620
$ws.Range("N5").Value = 'Post-conceptual artist exhibits a new selection of work at this Chelsea gallery'
$ws.Range("O5").Value = 'Best known for his word paintings—dropping vowels in pieces like "TRBL" and "DRNK"—the post-conceptual artist exhibits a new selection of work at this Chelsea gallery. Although Luhring Augustine has yet to announce details, we''d expect a muted palette and immediately discernable point of view. Wool''s paintings can be described in many ways, but shy isn''t one of them.'
$ws.Range("P5").Value = 'http://www.nycgo.com/images/460x285/Luhring-Augustine_V1_460x285.jpg'

$ws.Range("B2:B5").Select()
